$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计" (i.e.
#    before the current "2021-Q2" sheet), mirroring the workbook.xml
#    sheet-order change in the diff.
# ------------------------------------------------------------------
$q2sheet = $wb.Worksheets.Item("2021-Q2")
$newSheet = $wb.Worksheets.Add($q2sheet)
$newSheet.Name = "2022-Q4"

# The source data keeps fund codes / percentage-style figures as plain
# text (so leading/trailing zeros like "005457" or "2.80" survive).
# Force those ranges to Text format *before* writing the values so the
# engine doesn't silently coerce them to numbers.
$newSheet.Range("B2:B5").NumberFormat = "@"
$newSheet.Range("D2:F5").NumberFormat = "@"
$newSheet.Range("G2:G4").NumberFormat = "@"

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Row 2 - 005457
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "005457"
$newSheet.Cells.Item(2,3).Value = "景顺长城量化小盘股票"
$newSheet.Cells.Item(2,4).Value = "5.08"
$newSheet.Cells.Item(2,5).Value = "94.36"
$newSheet.Cells.Item(2,6).Value = "1.38"
$newSheet.Cells.Item(2,7).Value = "0.0701"
$newSheet.Cells.Item(2,8).Value = 6

# Row 3 - 202019
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "202019"
$newSheet.Cells.Item(3,3).Value = "南方策略优化混合"
$newSheet.Cells.Item(3,4).Value = "2.80"
$newSheet.Cells.Item(3,5).Value = "93.64"
$newSheet.Cells.Item(3,6).Value = "2.01"
$newSheet.Cells.Item(3,7).Value = "0.0563"
$newSheet.Cells.Item(3,8).Value = 9

# Row 4 - 014556
$newSheet.Cells.Item(4,1).Value = 2
$newSheet.Cells.Item(4,2).Value = "014556"
$newSheet.Cells.Item(4,3).Value = "富荣量化精选混合A"
$newSheet.Cells.Item(4,4).Value = "0.09"
$newSheet.Cells.Item(4,5).Value = "61.44"
$newSheet.Cells.Item(4,6).Value = "2.22"
$newSheet.Cells.Item(4,7).Value = "0.0020"
$newSheet.Cells.Item(4,8).Value = 10

# Row 5 - 014557
$newSheet.Cells.Item(5,1).Value = 3
$newSheet.Cells.Item(5,2).Value = "014557"
$newSheet.Cells.Item(5,3).Value = "富荣量化精选混合C"
$newSheet.Cells.Item(5,4).Value = "0.00"
$newSheet.Cells.Item(5,5).Value = "61.44"
$newSheet.Cells.Item(5,6).Value = "2.22"
$newSheet.Cells.Item(5,7).Value = 0
$newSheet.Cells.Item(5,8).Value = 10

# Header + "序号" index column styling: bold text, thin border all
# round, centered horizontally and top-aligned vertically (matches
# the look of the neighbouring quarter sheets). Apply cell-by-cell
# (rather than as one multi-cell range) so every cell gets the exact
# same uniform "thin box" border instead of an inside/outside border
# mix.
$styledCells = @()
for ($col = 2; $col -le 8; $col++) { $styledCells += $newSheet.Cells.Item(1, $col) }
for ($row = 2; $row -le 5; $row++) { $styledCells += $newSheet.Cells.Item($row, 1) }

foreach ($cell in $styledCells) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.Item(7).LineStyle = 1
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
    $cell.Borders.Item(10).LineStyle = 1
}

# ------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push existing data rows down by
#    one and insert the new 2022-Q4 summary row at row 2.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Capture existing rows 2 & 3 before overwriting anything.
$a2 = $total.Cells.Item(2,1).Value2
$b2 = $total.Cells.Item(2,2).Value2
$c2 = $total.Cells.Item(2,3).Value2
$d2 = $total.Cells.Item(2,4).Value2

$a3 = $total.Cells.Item(3,1).Value2
$b3 = $total.Cells.Item(3,2).Value2
$c3 = $total.Cells.Item(3,3).Value2
$d3 = $total.Cells.Item(3,4).Value2

# Move old row 2 ("2021-Q2") down to row 3.
$total.Cells.Item(3,1).Value = $a2
$total.Cells.Item(3,2).Value = $b2
$total.Cells.Item(3,3).Value = $c2
$total.Cells.Item(3,4).Value = $d2

# Move old row 3 ("2020-Q4") down to row 4.
$total.Cells.Item(4,1).Value = $a3
$total.Cells.Item(4,2).Value = $b3
$total.Cells.Item(4,3).Value = $c3
$total.Cells.Item(4,4).Value = $d3

# New row 2: 2022-Q4 summary.
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 0.13

# Row 4's "A" cell is brand-new (the sheet previously only had 3
# rows), so it has no inherited style yet; clone the "序号" column
# formatting from row 3 before setting its value.
$total.Cells.Item(3,1).Copy()
$total.Cells.Item(4,1).PasteSpecial(-4122)

# Fix up the "序号" (index) column for the rows that moved, matching
# the 0/1/2 sequence shown in the diff.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
